$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Akdis%Cezmi A%coreGivesNoEmail%0,                 Cao%Yi-yuan%coreGivesNoEmail%0,                 Dong%Xiang%coreGivesNoEmail%0,                 Gao%Ya-dong%coreGivesNoEmail%0,                 Yan%You-qin%coreGivesNoEmail%0,                 Yang%Yi-bin%coreGivesNoEmail%0,                 Yuan%Ya-dong%coreGivesNoEmail%0,                 Zhang%Jin-jin%coreGivesNoEmail%0]'
$ws.Range("E3").Value = '[Zhongliang%Wang%NULL%0,                   Bohan%Yang%NULL%0,                   Qianwen%Li%NULL%0,                   Lu%Wen%NULL%0,                   Ruiguang%Zhang%zrg27@163.com%0]'
$ws.Range("E4").Value = '[Xiaobo%Yang%NULL%0,                   Yuan%Yu%NULL%0,                   Jiqian%Xu%NULL%0,                   Huaqing%Shu%NULL%0,                   Jia''an%Xia%NULL%0,                   Hong%Liu%NULL%0,                   Yongran%Wu%NULL%0,                   Lu%Zhang%NULL%0,                   Zhui%Yu%NULL%0,                   Minghao%Fang%NULL%0,                   Ting%Yu%NULL%0,                   Yaxin%Wang%NULL%0,                   Shangwen%Pan%NULL%0,                   Xiaojing%Zou%NULL%0,                   Shiying%Yuan%NULL%0,                   You%Shang%NULL%0]'
$ws.Range("E5").Value = '[Chaomin%Wu%xref no email%0,           Xiaoyan%Chen%xref no email%0,           Yanping%Cai%xref no email%0,           Jia\u2019an%Xia%xref no email%0,           Xing%Zhou%xref no email%0,           Sha%Xu%xref no email%0,           Hanping%Huang%xref no email%0,           Li%Zhang%xref no email%0,           Xia%Zhou%xref no email%0,           Chunling%Du%xref no email%0,           Yuye%Zhang%xref no email%0,           Juan%Song%xref no email%0,           Sijiao%Wang%xref no email%0,           Yencheng%Chao%xref no email%0,           Zeyong%Yang%xref no email%0,           Jie%Xu%xref no email%0,           Xin%Zhou%xref no email%0,           Dechang%Chen%xref no email%0,           Weining%Xiong%xref no email%0,           Lei%Xu%xref no email%0,           Feng%Zhou%xref no email%0,           Jinjun%Jiang%xref no email%0,           Chunxue%Bai%xref no email%0,           Junhua%Zheng%xref no email%0,           Yuanlin%Song%xref no email%0]'
$ws.Range("E6").Value = '[Yihui%Huang%NULL%0,                   Mengqi%Tu%NULL%0,                   Shipei%Wang%NULL%0,                   Sichao%Chen%NULL%0,                   Wei%Zhou%NULL%0,                   Danyang%Chen%NULL%0,                   Lin%Zhou%NULL%0,                   Min%Wang%NULL%0,                   Yan%Zhao%NULL%0,                   Wen%Zeng%NULL%0,                   Qi%Huang%NULL%0,                   Hai''bo%Xu%NULL%0,                   Zeming%Liu%NULL%0,                   Liang%Guo%NULL%0]'
$ws.Range("E7").Value = '[Kunhua%Li%NULL%0,                   Jiong%Wu%NULL%0,                   Faqi%Wu%NULL%0,                   Dajing%Guo%NULL%0,                   Linli%Chen%NULL%0,                   Zheng%Fang%NULL%0,                   Chuanming%Li%NULL%0]'
$ws.Range("E8").Value = '[Xi%Xu%NULL%0,                   Chengcheng%Yu%NULL%0,                   Jing%Qu%NULL%0,                   Lieguang%Zhang%NULL%0,                   Songfeng%Jiang%NULL%0,                   Deyang%Huang%NULL%0,                   Bihua%Chen%NULL%0,                   Zhiping%Zhang%NULL%0,                   Wanhua%Guan%NULL%0,                   Zhoukun%Ling%NULL%0,                   Rui%Jiang%NULL%0,                   Tianli%Hu%NULL%0,                   Yan%Ding%NULL%0,                   Lin%Lin%NULL%0,                   Qingxin%Gan%NULL%0,                   Liangping%Luo%tluolp@jnu.edu.cn%0,                   Xiaoping%Tang%xtang@21cn.com%0,                   Jinxin%Liu%Liujx83710378@126.com%0]'
$ws.Range("E9").Value = '[Wenjie%Yang%NULL%0,                   Qiqi%Cao%NULL%0,                   Le%Qin%NULL%0,                   Xiaoyang%Wang%NULL%0,                   Zenghui%Cheng%NULL%0,                   Ashan%Pan%NULL%0,                   Jianyi%Dai%NULL%0,                   Qingfeng%Sun%NULL%0,                   Fengquan%Zhao%NULL%0,                   Jieming%Qu%NULL%0,                   Fuhua%Yan%NULL%0]'
$ws.Range("E10").Value = '[Wenjie%Yang%NULL%0,                   Qiqi%Cao%NULL%0,                   Le%Qin%NULL%0,                   Xiaoyang%Wang%NULL%0,                   Zenghui%Cheng%NULL%0,                   Ashan%Pan%NULL%0,                   Jianyi%Dai%NULL%0,                   Qingfeng%Sun%NULL%0,                   Fengquan%Zhao%NULL%0,                   Jieming%Qu%NULL%0,                   Fuhua%Yan%NULL%0]'
$ws.Range("E11").Value = '[Shaoqing%Lei%NULL%0,                   Fang%Jiang%NULL%0,                   Wating%Su%NULL%0,                   Chang%Chen%NULL%0,                   Jingli%Chen%NULL%0,                   Wei%Mei%NULL%0,                   Li-Ying%Zhan%NULL%0,                   Yifan%Jia%NULL%0,                   Liangqing%Zhang%NULL%0,                   Danyong%Liu%NULL%0,                   Zhong-Yuan%Xia%xiazhongyuan2005@aliyun.com%0,                   Zhengyuan%Xia%zyxia@hkucc.hku.hk%0]'
$ws.Range("E12").Value = '[Yun%Feng%NULL%0,                   Yun%Ling%NULL%0,                   Yun%Ling%NULL%0,                   Tao%Bai%NULL%0,                   Tao%Bai%NULL%0,                   Yusang%Xie%NULL%1,                   Jie%Huang%NULL%2,                   Jie%Huang%NULL%0,                   Jian%Li%NULL%1,                   Weining%Xiong%NULL%1,                   Dexiang%Yang%NULL%1,                   Rong%Chen%NULL%1,                   Fangying%Lu%NULL%1,                   Yunfei%Lu%NULL%2,                   Xuhui%Liu%NULL%1,                   Yuqing%Chen%NULL%2,                   Yuqing%Chen%NULL%0,                   Xin%Li%NULL%0,                   Yong%Li%NULL%1,                   Hanssa Dwarka%Summah%NULL%1,                   Huihuang%Lin%NULL%1,                   Jiayang%Yan%NULL%1,                   Min%Zhou%NULL%0,                   Hongzhou%Lu%NULL%0,                   Hongzhou%Lu%NULL%0,                   Jieming%Qu%NULL%0,                   Jieming%Qu%NULL%0]'
$ws.Range("E13").Value = '[Mingli%Yuan%NULL%1,                   Wen%Yin%NULL%0,                   Wen%Yin%NULL%0,                   Zhaowu%Tao%NULL%1,                   Weijun%Tan%NULL%1,                   Yi%Hu%NULL%0,                   Oliver%Schildgen%NULL%2,                   Oliver%Schildgen%NULL%0]'
$ws.Range("E14").Value = '[Pingzheng%Mo%NULL%0,                   Yuanyuan%Xing%NULL%0,                   Yu%Xiao%NULL%0,                   Liping%Deng%NULL%0,                   Qiu%Zhao%NULL%0,                   Hongling%Wang%NULL%0,                   Yong%Xiong%NULL%0,                   Zhenshun%Cheng%NULL%0,                   Shicheng%Gao%NULL%0,                   Ke%Liang%NULL%0,                   Mingqi%Luo%NULL%0,                   Tielong%Chen%NULL%0,                   Shihui%Song%NULL%0,                   Zhiyong%Ma%NULL%0,                   Xiaoping%Chen%NULL%0,                   Ruiying%Zheng%NULL%0,                   Qian%Cao%NULL%0,                   Fan%Wang%fanndywang@foxmail.com%0,                   Yongxi%Zhang%znact1936@126.com%0]'
$ws.Range("E15").Value = '[Luwen%Wang%NULL%0,                   Xun%Li%NULL%1,                   Hui%Chen%NULL%0,                   Shaonan%Yan%NULL%1,                   Dong%Li%NULL%1,                   Yan%Li%NULL%1,                   Zuojiong%Gong%NULL%1]'
$ws.Range("E16").Value = '[Guqin%Zhang%NULL%0,                   Chang%Hu%NULL%0,                   Linjie%Luo%NULL%0,                   Fang%Fang%NULL%0,                   Yongfeng%Chen%NULL%0,                   Jianguo%Li%NULL%0,                   Zhiyong%Peng%NULL%0,                   Huaqin%Pan%phq2012@whu.edu.cn%0]'
$ws.Range("E17").Value = '[Tao%Guo%NULL%0,                   Yongzhen%Fan%NULL%0,                   Ming%Chen%NULL%0,                   Xiaoyan%Wu%NULL%0,                   Lin%Zhang%NULL%0,                   Tao%He%NULL%0,                   Hairong%Wang%NULL%0,                   Jing%Wan%NULL%0,                   Xinghuan%Wang%NULL%0,                   Zhibing%Lu%NULL%0]'
$ws.Range("E18").Value = '[Barnaby%D.P.%coreGivesNoEmail%0,                 Becker%L.B.%coreGivesNoEmail%0,                 Chelico%J.D.%coreGivesNoEmail%0,                 Cohen%S.L.%coreGivesNoEmail%0,                 Cookingham%J.%coreGivesNoEmail%0,                 Coppa%K.%coreGivesNoEmail%0,                 Crawford%J.M.%coreGivesNoEmail%0,                 Davidson%K.W.%coreGivesNoEmail%0,                 Diefenbach%M.A.%coreGivesNoEmail%0,                 Dominello%A.J.%coreGivesNoEmail%0,                 Duer-Hefele%J.%coreGivesNoEmail%0,                 Falzon%L.%coreGivesNoEmail%0,                 Gitlin%J.%coreGivesNoEmail%0,                 Hajizadeh%N.%coreGivesNoEmail%0,                 Harvin%T.G.%coreGivesNoEmail%0,                 Hirsch%J.S.%coreGivesNoEmail%0,                 Hirschwerk%D.A.%coreGivesNoEmail%0,                 Kim%E.J.%coreGivesNoEmail%0,                 Kozel%Z.M.%coreGivesNoEmail%0,                 Marrast%L.M.%coreGivesNoEmail%0,                 McGinn%T.%coreGivesNoEmail%0,                 Mogavero%J.N.%coreGivesNoEmail%0,                 Narasimhan%M.%coreGivesNoEmail%0,                 Osorio%G.A.%coreGivesNoEmail%0,                 Qiu%M.%coreGivesNoEmail%0,                 Richardson%S.%coreGivesNoEmail%0,                 Zanos%T.P.%coreGivesNoEmail%0]'
$ws.Range("E19").Value = '[NULL%NULL%NULL%0,                   NULL%NULL%NULL%0,                   NULL%NULL%NULL%0,                   Nancy%Chow%NULL%0,                   Katherine%Fleming-Dutra%NULL%1,                   Ryan%Gierke%NULL%0,                   Aron%Hall%NULL%0,                   Michelle%Hughes%NULL%1,                   Tamara%Pilishvili%NULL%0,                   Matthew%Ritchey%NULL%0,                   Katherine%Roguski%NULL%1,                   Tami%Skoff%NULL%1,                   Emily%Ussery%NULL%1]'
$ws.Range("E20").Value = '[Barnaby Edward%Young%xref no email%0,           Sean Wei Xiang%Ong%xref no email%0,           Shirin%Kalimuddin%xref no email%0,           Jenny G.%Low%xref no email%0,           Seow Yen%Tan%xref no email%0,           Jiashen%Loh%xref no email%0,           Oon-Tek%Ng%xref no email%0,           Kalisvar%Marimuthu%xref no email%0,           Li Wei%Ang%xref no email%0,           Tze Minn%Mak%xref no email%0,           Sok Kiang%Lau%xref no email%0,           Danielle E.%Anderson%xref no email%0,           Kian Sing%Chan%xref no email%0,           Thean Yen%Tan%xref no email%0,           Tong Yong%Ng%xref no email%0,           Lin%Cui%xref no email%0,           Zubaidah%Said%xref no email%0,           Lalitha%Kurupatham%xref no email%0,           Mark I-Cheng%Chen%xref no email%0,           Monica%Chan%xref no email%0,           Shawn%Vasoo%xref no email%0,           Lin-Fa%Wang%xref no email%0,           Boon Huan%Tan%xref no email%0,           Raymond Tzer Pin%Lin%xref no email%0,           Vernon Jian Ming%Lee%xref no email%0,           Yee-Sin%Leo%xref no email%0,           David Chien%Lye%xref no email%0]'
$ws.Range("E22").Value = '[Antonelli%Massimo%coreGivesNoEmail%0,                 Cabrini%Luca%coreGivesNoEmail%0,                 Castelli%Antonio%coreGivesNoEmail%0,                 Cecconi%Maurizio%coreGivesNoEmail%0,                 Cereda%Danilo%coreGivesNoEmail%0,                 Coluccello%Antonio%coreGivesNoEmail%0,                 Foti%Giuseppe%coreGivesNoEmail%0,                 Fumagalli%Roberto%coreGivesNoEmail%0,                 Grasselli%Giacomo%coreGivesNoEmail%0,                 Iotti%Giorgio%coreGivesNoEmail%0,                 Latronico%Nicola%coreGivesNoEmail%0,                 Lorini%Luca%coreGivesNoEmail%0,                 Merler%Stefano%coreGivesNoEmail%0,                 Natalini%Giuseppe%coreGivesNoEmail%0,                 Pesenti%Antonio%coreGivesNoEmail%0,                 Piatti%Alessandra%coreGivesNoEmail%0,                 Ranieri%Marco Vito%coreGivesNoEmail%0,                 Scandroglio%Anna Mara%coreGivesNoEmail%0,                 Storti%Enrico%coreGivesNoEmail%0,                 Zanella%Alberto%coreGivesNoEmail%0,                 Zangrillo%Alberto%coreGivesNoEmail%0]'
$ws.Range("E23").Value = '[Yingzhen%Du%NULL%0,                   Lei%Tu%NULL%0,                   Lei%Tu%NULL%0,                   Pingjun%Zhu%NULL%0,                   Pingjun%Zhu%NULL%0,                   Mi%Mu%NULL%0,                   Mi%Mu%NULL%0,                   Runsheng%Wang%NULL%0,                   Pengcheng%Yang%NULL%0,                   Xi%Wang%NULL%0,                   Chao%Hu%NULL%0,                   Rongyu%Ping%NULL%0,                   Peng%Hu%NULL%0,                   Tianzhi%Li%NULL%0,                   Feng%Cao%NULL%0,                   Christopher%Chang%NULL%0,                   Qinyong%Hu%NULL%0,                   Yang%Jin%NULL%0,                   Guogang%Xu%NULL%0]'
$ws.Range("E25").Value = '[Nanshan%Chen%NULL%0,                   Min%Zhou%NULL%0,                   Xuan%Dong%NULL%0,                   Jieming%Qu%NULL%0,                   Fengyun%Gong%NULL%0,                   Yang%Han%NULL%0,                   Yang%Qiu%NULL%0,                   Jingli%Wang%NULL%0,                   Ying%Liu%NULL%0,                   Yuan%Wei%NULL%0,                   Jia''an%Xia%NULL%0,                   Ting%Yu%NULL%0,                   Xinxin%Zhang%NULL%0,                   Li%Zhang%NULL%0]'
$ws.Range("E26").Value = '[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,                 Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%2,                 Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%2,                 Subert-Salas%Lizandra%coreGivesNoEmail%2]'
$ws.Range("E28").Value = '[Chaolin%Huang%NULL%0,                   Yeming%Wang%NULL%0,                   Xingwang%Li%NULL%0,                   Lili%Ren%NULL%0,                   Jianping%Zhao%NULL%0,                   Yi%Hu%NULL%0,                   Li%Zhang%NULL%0,                   Guohui%Fan%NULL%0,                   Jiuyang%Xu%NULL%0,                   Xiaoying%Gu%NULL%0,                   Zhenshun%Cheng%NULL%0,                   Ting%Yu%NULL%0,                   Jiaan%Xia%NULL%0,                   Yuan%Wei%NULL%0,                   Wenjuan%Wu%NULL%0,                   Xuelei%Xie%NULL%0,                   Wen%Yin%NULL%0,                   Hui%Li%NULL%0,                   Min%Liu%NULL%0,                   Yan%Xiao%NULL%0,                   Hong%Gao%NULL%0,                   Li%Guo%NULL%0,                   Jungang%Xie%NULL%0,                   Guangfa%Wang%NULL%0,                   Rongmeng%Jiang%NULL%0,                   Zhancheng%Gao%NULL%0,                   Qi%Jin%NULL%0,                   Jianwei%Wang%wangjw28@163.com%0,                   Bin%Cao%caobin_ben@163.com%0]'
$ws.Range("E29").Value = '[Kui%Liu%NULL%0,                   Yuan-Yuan%Fang%NULL%0,                   Yan%Deng%NULL%0,                   Wei%Liu%NULL%0,                   Mei-Fang%Wang%NULL%0,                   Jing-Ping%Ma%NULL%0,                   Wei%Xiao%NULL%0,                   Ying-Nan%Wang%NULL%0,                   Min-Hua%Zhong%NULL%0,                   Cheng-Hong%Li%NULL%0,                   Guang-Cai%Li%NULL%0,                   Hui-Guo%Liu%NULL%0,                   Xiu-Yuan%Hao%NULL%0,                   Xiu-Yuan%Hao%NULL%0,                   Pei-Fang%Wei%NULL%0]'
$ws.Range("E30").Value = '[Yingxia%Liu%yingxialiu@hotmail.com%0,                   Yang%Yang%NULL%0,                   Cong%Zhang%NULL%0,                   Fengming%Huang%NULL%0,                   Fuxiang%Wang%NULL%0,                   Jing%Yuan%NULL%0,                   Zhaoqin%Wang%NULL%0,                   Jinxiu%Li%NULL%0,                   Jianming%Li%NULL%0,                   Cheng%Feng%NULL%0,                   Zheng%Zhang%NULL%0,                   Lifei%Wang%NULL%0,                   Ling%Peng%NULL%0,                   Li%Chen%NULL%0,                   Yuhao%Qin%NULL%0,                   Dandan%Zhao%NULL%0,                   Shuguang%Tan%NULL%0,                   Lu%Yin%NULL%0,                   Jun%Xu%NULL%0,                   Congzhao%Zhou%NULL%0,                   Chengyu%Jiang%jiang@pumc.edu.cn%0,                   Lei%Liu%liulei3322@aliyun.com%0]'
$ws.Range("E31").Value = '[Wei-jie%Guan%NULL%0,                   Zheng-yi%Ni%NULL%0,                   Zheng-yi%Ni%NULL%0,                   Yu%Hu%NULL%0,                   Wen-hua%Liang%NULL%0,                   Chun-quan%Ou%NULL%0,                   Jian-xing%He%NULL%0,                   Lei%Liu%NULL%0,                   Hong%Shan%NULL%0,                   Chun-liang%Lei%NULL%0,                   David S.C.%Hui%NULL%0,                   Bin%Du%NULL%0,                   Lan-juan%Li%NULL%0,                   Guang%Zeng%NULL%0,                   Kwok-Yung%Yuen%NULL%0,                   Ru-chong%Chen%NULL%0,                   Chun-li%Tang%NULL%0,                   Tao%Wang%NULL%0,                   Ping-yan%Chen%NULL%0,                   Jie%Xiang%NULL%0,                   Shi-yue%Li%NULL%0,                   Jin-lin%Wang%NULL%0,                   Zi-jing%Liang%NULL%0,                   Yi-xiang%Peng%NULL%0,                   Li%Wei%NULL%0,                   Yong%Liu%NULL%0,                   Ya-hua%Hu%NULL%0,                   Peng%Peng%NULL%0,                   Jian-ming%Wang%NULL%0,                   Ji-yang%Liu%NULL%0,                   Zhong%Chen%NULL%0,                   Gang%Li%NULL%0,                   Zhi-jian%Zheng%NULL%0,                   Shao-qin%Qiu%NULL%0,                   Jie%Luo%NULL%0,                   Chang-jiang%Ye%NULL%0,                   Shao-yong%Zhu%NULL%0,                   Nan-shan%Zhong%NULL%0]'
$ws.Range("E32").Value = '[Rong-Hui%Du%NULL%0,                   Li-Min%Liu%NULL%0,                   Wen%Yin%NULL%0,                   Wen%Wang%NULL%0,                   Lu-Lu%Guan%NULL%0,                   Ming-Li%Yuan%NULL%0,                   Yu-Lei%Li%NULL%0,                   Yi%Hu%NULL%0,                   Xu-Yan%Li%NULL%0,                   Bing%Sun%NULL%0,                   Peng%Peng%NULL%0,                   Huan-Zhong%Shi%NULL%0]'
$ws.Range("E33").Value = '[Wei-jie%Guan%NULL%0,                   Wen-hua%Liang%NULL%0,                   Yi%Zhao%NULL%0,                   Heng-rui%Liang%NULL%0,                   Zi-sheng%Chen%NULL%0,                   Yi-min%Li%NULL%0,                   Xiao-qing%Liu%NULL%0,                   Ru-chong%Chen%NULL%0,                   Chun-li%Tang%NULL%0,                   Tao%Wang%NULL%0,                   Chun-quan%Ou%NULL%0,                   Li%Li%NULL%0,                   Ping-yan%Chen%NULL%0,                   Ling%Sang%NULL%0,                   Wei%Wang%NULL%0,                   Jian-fu%Li%NULL%0,                   Cai-chen%Li%NULL%0,                   Li-min%Ou%NULL%0,                   Bo%Cheng%NULL%0,                   Shan%Xiong%NULL%0,                   Zheng-yi%Ni%NULL%0,                   Jie%Xiang%NULL%0,                   Yu%Hu%NULL%0,                   Lei%Liu%NULL%0,                   Hong%Shan%NULL%0,                   Chun-liang%Lei%NULL%0,                   Yi-xiang%Peng%NULL%0,                   Li%Wei%NULL%0,                   Yong%Liu%NULL%0,                   Ya-hua%Hu%NULL%0,                   Peng%Peng%NULL%0,                   Jian-ming%Wang%NULL%0,                   Ji-yang%Liu%NULL%0,                   Zhong%Chen%NULL%0,                   Gang%Li%NULL%0,                   Zhi-jian%Zheng%NULL%0,                   Shao-qin%Qiu%NULL%0,                   Jie%Luo%NULL%0,                   Chang-jiang%Ye%NULL%0,                   Shao-yong%Zhu%NULL%0,                   Lin-ling%Cheng%NULL%0,                   Feng%Ye%NULL%0,                   Shi-yue%Li%NULL%0,                   Jin-ping%Zheng%NULL%0,                   Nuo-fu%Zhang%NULL%0,                   Nan-shan%Zhong%NULL%0,                   Jian-xing%He%NULL%0]'
$ws.Range("E34").Value = '[Xiao-Wei%Xu%NULL%0,                   Xiao-Xin%Wu%NULL%0,                   Xian-Gao%Jiang%NULL%0,                   Kai-Jin%Xu%NULL%0,                   Ling-Jun%Ying%NULL%0,                   Chun-Lian%Ma%NULL%0,                   Shi-Bo%Li%NULL%0,                   Hua-Ying%Wang%NULL%0,                   Sheng%Zhang%NULL%0,                   Hai-Nv%Gao%NULL%0,                   Ji-Fang%Sheng%NULL%0,                   Hong-Liu%Cai%NULL%0,                   Yun-Qing%Qiu%NULL%0,                   Lan-Juan%Li%NULL%0]'
$ws.Range("E35").Value = '[TieLong%Chen%NULL%0,                   Zhe%Dai%NULL%1,                   Pingzheng%Mo%NULL%1,                   Xinyu%Li%NULL%1,                   Zhiyong%Ma%NULL%0,                   Shihui%Song%NULL%0,                   Xiaoping%Chen%NULL%0,                   Mingqi%Luo%NULL%0,                   Ke%Liang%NULL%0,                   Shicheng%Gao%NULL%0,                   Yongxi%Zhang%NULL%0,                   Liping%Deng%dengdeng78@126.com%0,                   Yong%Xiong%NULL%0,                   Yong%Xiong%NULL%0]'
$ws.Range("E36").Value = '[Wen-hua%Liang%NULL%0,                   Wei-jie%Guan%NULL%0,                   Cai-chen%Li%NULL%0,                   Yi-min%Li%NULL%0,                   Heng-rui%Liang%NULL%0,                   Yi%Zhao%NULL%0,                   Xiao-qing%Liu%NULL%0,                   Ling%Sang%NULL%0,                   Ru-chong%Chen%NULL%0,                   Chun-li%Tang%NULL%0,                   Tao%Wang%NULL%0,                   Wei%Wang%NULL%0,                   Qi-hua%He%NULL%0,                   Zi-sheng%Chen%NULL%0,                   Sook-San%Wong%NULL%0,                   Mark%Zanin%NULL%0,                   Jun%Liu%NULL%0,                   Xin%Xu%NULL%0,                   Jun%Huang%NULL%0,                   Jian-fu%Li%NULL%0,                   Li-min%Ou%NULL%0,                   Bo%Cheng%NULL%0,                   Shan%Xiong%NULL%0,                   Zhan-hong%Xie%NULL%0,                   Zheng-yi%Ni%NULL%0,                   Yu%Hu%NULL%0,                   Lei%Liu%NULL%0,                   Hong%Shan%NULL%0,                   Chun-liang%Lei%NULL%0,                   Yi-xiang%Peng%NULL%0,                   Li%Wei%NULL%0,                   Yong%Liu%NULL%0,                   Ya-hua%Hu%NULL%0,                   Peng%Peng%NULL%0,                   Jian-ming%Wang%NULL%0,                   Ji-yang%Liu%NULL%0,                   Zhong%Chen%NULL%0,                   Gang%Li%NULL%0,                   Zhi-jian%Zheng%NULL%0,                   Shao-qin%Qiu%NULL%0,                   Jie%Luo%NULL%0,                   Chang-jiang%Ye%NULL%0,                   Shao-yong%Zhu%NULL%0,                   Lin-ling%Cheng%NULL%0,                   Feng%Ye%NULL%0,                   Shi-yue%Li%NULL%0,                   Jin-ping%Zheng%NULL%0,                   Nuo-fu%Zhang%NULL%0,                   Nan-shan%Zhong%NULL%0,                   Jian-xing%He%NULL%0]'
$ws.Range("E37").Value = '[Lang%Wang%NULL%0,                   Wenbo%He%NULL%0,                   Xiaomei%Yu%NULL%0,                   Dalong%Hu%NULL%0,                   Mingwei%Bao%NULL%0,                   Huafen%Liu%NULL%0,                   Jiali%Zhou%NULL%0,                   Hong%Jiang%NULL%0]'
$ws.Range("E38").Value = '[Jianlei%Cao%NULL%0,                   Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%0,                   Wenlin%Cheng%NULL%0,                   Lei%Yu%NULL%0,                   Ya-Kun%Liu%NULL%0,                   Xiaoyong%Hu%NULL%0,                   Qiang%Liu%NULL%0]'
$ws.Range("E39").Value = '[Fei%Zhou%NULL%0,                   Ting%Yu%NULL%0,                   Ronghui%Du%NULL%0,                   Guohui%Fan%NULL%0,                   Ying%Liu%NULL%0,                   Zhibo%Liu%NULL%0,                   Jie%Xiang%NULL%0,                   Yeming%Wang%NULL%0,                   Bin%Song%NULL%0,                   Xiaoying%Gu%NULL%0,                   Lulu%Guan%NULL%0,                   Yuan%Wei%NULL%0,                   Hui%Li%NULL%0,                   Xudong%Wu%NULL%0,                   Jiuyang%Xu%NULL%0,                   Shengjin%Tu%NULL%0,                   Yi%Zhang%NULL%0,                   Hua%Chen%NULL%0,                   Bin%Cao%NULL%0]'
$ws.Range("E40").Value = '[Rui%Zhang%NULL%0,                   Huangqing%Ouyang%NULL%1,                   Lingli%Fu%NULL%1,                   Shijie%Wang%NULL%1,                   Jianglong%Han%NULL%1,                   Kejie%Huang%NULL%1,                   Mingfang%Jia%NULL%1,                   Qibin%Song%NULL%1,                   Zhenming%Fu%davidfuzming@whu.edu.cn%1]'
$ws.Range("E41").Value = '[Tao%Chen%NULL%0,                   Di%Wu%NULL%0,                   Huilong%Chen%NULL%0,                   Weiming%Yan%NULL%0,                   Danlei%Yang%NULL%0,                   Guang%Chen%NULL%0,                   Ke%Ma%NULL%0,                   Dong%Xu%NULL%0,                   Haijing%Yu%NULL%0,                   Hongwu%Wang%NULL%0,                   Tao%Wang%NULL%0,                   Wei%Guo%NULL%0,                   Jia%Chen%NULL%0,                   Chen%Ding%NULL%0,                   Xiaoping%Zhang%NULL%0,                   Jiaquan%Huang%NULL%0,                   Meifang%Han%NULL%0,                   Shusheng%Li%NULL%0,                   Xiaoping%Luo%NULL%0,                   Jianping%Zhao%NULL%0,                   Qin%Ning%NULL%0]'
